$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new timestamp text (column B) and new value (column C, $null = unchanged)
$updates = @(
    @{ Row = 2;  Time = "2022-03-15 15:01:56.185909"; Value = 99 },
    @{ Row = 3;  Time = "2022-03-15 15:01:57.107958"; Value = $null },
    @{ Row = 4;  Time = "2022-03-15 15:01:58.700224"; Value = 88 },
    @{ Row = 5;  Time = "2022-03-15 15:02:01.189186"; Value = 80 },
    @{ Row = 6;  Time = "2022-03-15 15:02:02.684583"; Value = 77 },
    @{ Row = 7;  Time = "2022-03-15 15:02:03.720169"; Value = 70 },
    @{ Row = 8;  Time = "2022-03-15 15:02:04.413693"; Value = 66 },
    @{ Row = 9;  Time = "2022-03-15 15:02:06.472708"; Value = 60 },
    @{ Row = 10; Time = "2022-03-15 15:02:07.204546"; Value = 55 },
    @{ Row = 11; Time = "2022-03-15 15:02:07.953563"; Value = 50 },
    @{ Row = 12; Time = "2022-03-15 15:02:08.667783"; Value = 44 },
    @{ Row = 13; Time = "2022-03-15 15:02:09.492311"; Value = 40 },
    @{ Row = 14; Time = "2022-03-15 15:02:10.457625"; Value = 33 },
    @{ Row = 15; Time = "2022-03-15 15:02:11.600791"; Value = 30 },
    @{ Row = 16; Time = "2022-03-15 15:02:12.681954"; Value = 22 },
    @{ Row = 17; Time = "2022-03-15 15:02:13.545104"; Value = 20 },
    @{ Row = 18; Time = "2022-03-15 15:02:14.608945"; Value = 11 },
    @{ Row = 19; Time = "2022-03-15 15:02:15.523652"; Value = 10 },
    @{ Row = 20; Time = "2022-03-15 15:02:16.385645"; Value = 9 },
    @{ Row = 21; Time = "2022-03-15 15:02:16.842992"; Value = 8 },
    @{ Row = 22; Time = "2022-03-15 15:02:17.298025"; Value = 7 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.Time
    if ($null -ne $u.Value) {
        $ws.Cells.Item($r, 3).Value = $u.Value
    }
}
